$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 298
$ws.Range("I53").Value = 196.57143
$ws.Range("J53").Value = 386.75
$ws.Range("K53").Value = 196.57143
$ws.Range("L53").Value = 386.75
$ws.Range("M53").Value = 440.42857
$ws.Range("N53").Value = -1660.75
$ws.Range("H133").Value = 12300
$ws.Range("J133").Value = 12300
$ws.Range("L133").Value = 12300
$ws.Range("N133").Value = -22420

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5080.3213
$ws.Range("J2").Value = 1060.1428
$ws.Range("L2").Value = 1060.1428
$ws.Range("N2").Value = -1286.1428
$ws.Range("H3").Value = 2420
$ws.Range("J3").Value = 2966.6667
$ws.Range("L3").Value = 2966.6667
$ws.Range("N3").Value = -3196.6667
$ws.Range("H5").Value = 1667025
$ws.Range("I5").Value = 1667025
$ws.Range("K5").Value = 1667025
$ws.Range("M5").Value = -1666913
$ws.Range("H32").Value = 2315.353
$ws.Range("I32").Value = 1841.591
$ws.Range("K32").Value = 1841.591
$ws.Range("M32").Value = -1554.591
$ws.Range("H45").Value = 2259.4707
$ws.Range("I45").Value = 1250.9166
$ws.Range("J45").Value = 4680
$ws.Range("K45").Value = 1250.9166
$ws.Range("L45").Value = 4680
$ws.Range("M45").Value = -873.9166
$ws.Range("N45").Value = -5434
$ws.Range("H116").Value = 5080.3213
$ws.Range("J116").Value = 1060.1428
$ws.Range("L116").Value = 1060.1428
$ws.Range("N116").Value = -5648.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5080.3213
$ws.Range("J3").Value = 1060.1428
$ws.Range("L3").Value = 1060.1428
$ws.Range("N3").Value = -1288.1428
$ws.Range("H4").Value = 1667025
$ws.Range("I4").Value = 1667025
$ws.Range("K4").Value = 1667025
$ws.Range("M4").Value = -1666910
$ws.Range("H59").Value = 34900
$ws.Range("J59").Value = 34900
$ws.Range("L59").Value = 34900
$ws.Range("N59").Value = -36594
$ws.Range("H76").Value = 49000
$ws.Range("J76").Value = 49000
$ws.Range("L76").Value = 49000
$ws.Range("N76").Value = -49630
$ws.Range("H79").Value = 49000
$ws.Range("J79").Value = 49000
$ws.Range("L79").Value = 49000
$ws.Range("N79").Value = -51184
$ws.Range("H99").Value = 3027.125
$ws.Range("I99").Value = 525.4
$ws.Range("J99").Value = 7196.6665
$ws.Range("K99").Value = 525.4
$ws.Range("L99").Value = 7196.6665
$ws.Range("M99").Value = 972.6
$ws.Range("N99").Value = -10192.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 277.16666
$ws.Range("I22").Value = 277.16666
$ws.Range("K22").Value = 277.16666
$ws.Range("M22").Value = 72.83334000000002
$ws.Range("H31").Value = 2039.7354
$ws.Range("I31").Value = 1268.7693
$ws.Range("K31").Value = 1268.7693
$ws.Range("M31").Value = -973.7692999999999
$ws.Range("H34").Value = 2039.7354
$ws.Range("I34").Value = 1268.7693
$ws.Range("K34").Value = 1268.7693
$ws.Range("M34").Value = -1066.7693
$ws.Range("H58").Value = 933.4151000000001
$ws.Range("I58").Value = 490.44736
$ws.Range("J58").Value = 2055.6
$ws.Range("K58").Value = 490.44736
$ws.Range("L58").Value = 2055.6
$ws.Range("M58").Value = -287.44736
$ws.Range("N58").Value = -2461.6
$ws.Range("H94").Value = 1087.5416
$ws.Range("I94").Value = 1039.4
$ws.Range("J94").Value = 1100.2106
$ws.Range("K94").Value = 1039.4
$ws.Range("L94").Value = 1100.2106
$ws.Range("M94").Value = -588.4000000000001
$ws.Range("N94").Value = -2002.2106
$ws.Range("H136").Value = 933.4151000000001
$ws.Range("I136").Value = 490.44736
$ws.Range("J136").Value = 2055.6
$ws.Range("K136").Value = 1471.34208
$ws.Range("L136").Value = 6166.799999999999
$ws.Range("M136").Value = 1078.65792
$ws.Range("N136").Value = -11266.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 999
$ws.Range("J98").Value = 999
$ws.Range("L98").Value = 2997
$ws.Range("N98").Value = -5993
$ws.Range("H131").Value = 2611.7888
$ws.Range("I131").Value = 348.8889
$ws.Range("J131").Value = 2940.2742
$ws.Range("K131").Value = 1046.6667
$ws.Range("L131").Value = 8820.8226
$ws.Range("M131").Value = 3993.3333
$ws.Range("N131").Value = -18900.8226

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 45
$ws.Range("I2").Value = 30
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 30
$ws.Range("L2").Value = 60
$ws.Range("M2").Value = 83
$ws.Range("N2").Value = -286
$ws.Range("H126").Value = 2472.2222
$ws.Range("I126").Value = 2400
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 7200
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -4730
$ws.Range("N126").Value = -12440
$ws.Range("H137").Value = 50500
$ws.Range("J137").Value = 50500
$ws.Range("L137").Value = 50500
$ws.Range("N137").Value = -60700
$ws.Range("H139").Value = 59515.6
$ws.Range("J139").Value = 59515.6
$ws.Range("L139").Value = 59515.6
$ws.Range("N139").Value = -69795.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 50000
$ws.Range("I18").Value = 50000
$ws.Range("K18").Value = 50000
$ws.Range("M18").Value = -49828
$ws.Range("H22").Value = 1250
$ws.Range("I22").Value = 1333.3334
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1333.3334
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -1038.3334
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 1250
$ws.Range("I27").Value = 1333.3334
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 1333.3334
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -1226.3334
$ws.Range("N27").Value = -1214
$ws.Range("H40").Value = 4537.8823
$ws.Range("I40").Value = 3381.3333
$ws.Range("K40").Value = 3381.3333
$ws.Range("M40").Value = -3245.3333
$ws.Range("H46").Value = 1024.8572
$ws.Range("I46").Value = 425
$ws.Range("J46").Value = 1124.8334
$ws.Range("K46").Value = 425
$ws.Range("L46").Value = 1124.8334
$ws.Range("M46").Value = -237
$ws.Range("N46").Value = -1500.8334
$ws.Range("H48").Value = 22497.25
$ws.Range("I48").Value = 26500
$ws.Range("J48").Value = 18494.5
$ws.Range("K48").Value = 26500
$ws.Range("L48").Value = 18494.5
$ws.Range("M48").Value = -25839
$ws.Range("N48").Value = -19816.5
$ws.Range("H100").Value = 2786.8948
$ws.Range("J100").Value = 2775.0557
$ws.Range("L100").Value = 2775.0557
$ws.Range("N100").Value = -3857.0557
$ws.Range("H132").Value = 4293.1943
$ws.Range("J132").Value = 5189.3335
$ws.Range("L132").Value = 15568.0005
$ws.Range("N132").Value = -20628.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 27948.25
$ws.Range("I62").Value = 51500
$ws.Range("J62").Value = 20097.666
$ws.Range("K62").Value = 51500
$ws.Range("L62").Value = 20097.666
$ws.Range("M62").Value = -50876
$ws.Range("N62").Value = -21345.666
$ws.Range("H65").Value = 27948.25
$ws.Range("I65").Value = 51500
$ws.Range("J65").Value = 20097.666
$ws.Range("K65").Value = 257500
$ws.Range("L65").Value = 100488.33
$ws.Range("M65").Value = -254380
$ws.Range("N65").Value = -106728.33
$ws.Range("H81").Value = 62333.47
$ws.Range("I81").Value = 168498.17
$ws.Range("J81").Value = 4425.4546
$ws.Range("K81").Value = 336996.34
$ws.Range("L81").Value = 8850.9092
$ws.Range("M81").Value = -335935.34
$ws.Range("N81").Value = -10972.9092
$ws.Range("H84").Value = 62333.47
$ws.Range("I84").Value = 168498.17
$ws.Range("J84").Value = 4425.4546
$ws.Range("K84").Value = 1684981.7
$ws.Range("L84").Value = 44254.546
$ws.Range("M84").Value = -1679677.7
$ws.Range("N84").Value = -54862.546
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 63333.332
$ws.Range("J139").Value = 63333.332
$ws.Range("L139").Value = 63333.332
$ws.Range("N139").Value = -73613.33199999999
